$d = $word.ActiveDocument

# Locate the sentence being rewritten and expand the range to the
# whole paragraph that contains it.
$old = "which produces a single value that changes no matter which axis the IMU is rotated around"
$rng = $d.Content
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the sentence to rewrite"
}
[void]$rng.Expand(4)

# Exclude the trailing paragraph mark so InsertXML only rewrites the
# paragraph's run content, preserving the <w:p> element itself.
[void]$rng.MoveEnd(1, -1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">The BNO055 is run using Adafruit’s sensor libraries. An I2C communication line receives the measurements from the IMU, and the values are read in the form of a Quaternion, </w:t></w:r><w:r><w:t>a vector type that combines the rotation values around the Cartesian axes to</w:t></w:r><w:r><w:t xml:space="preserve"> produce a single value </w:t></w:r><w:r><w:t>representing the rotation around an arbitrary axis</w:t></w:r><w:r><w:t xml:space="preserve"> known as the Euler axis</w:t></w:r><w:r><w:t>. This makes it easier to compare the data against a rolling average of previous readings, allowing the device to not produce false positives from small bumps to the surface the device is on, while also making it harder to move the device without setting off the alarm.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$rng.InsertXML($xml)
"Quaternion explanation expanded successfully"
